$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add a new bullet paragraph (ListParagraph / numId=4) right after "The
#    centre/origin of frame is at the centre of the phone." and before the
#    "Mechanical and 3D printing" Heading1.
# ---------------------------------------------------------------------------
$anchorRng = $d.Content
$anchorRng.Find.Execute("The centre/origin of frame is at the centre of the phone.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorRng.Collapse(0)          # wdCollapseEnd
$anchorRng.InsertParagraphAfter()

# Locate the freshly-minted (still empty) paragraph; it inherited the
# ListParagraph / numId=4 formatting from the paragraph above it.
$newPara = $anchorRng.Next(4, 1).Paragraphs.First
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq $anchorRng.End) {
        $newPara = $cand
        break
    }
}

$newRng = $newPara.Range
$newRng.MoveEnd(1, -1)          # wdCharacter, exclude the paragraph mark

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">See the folder example_code/from_zac_apelt, which has the instructions for perspective transform. </w:t></w:r><w:r><w:t xml:space="preserve">Read the .jpg picture first. </w:t></w:r><w:r><w:t>Zac did all of it last year, which is why I don' + [char]0x2019 + 't know much about perspective transform.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Mark "3D print must be RED..." as now starting a rendered page (insert
#    <w:lastRenderedPageBreak/> as the first child of that run).
# ---------------------------------------------------------------------------
$rngRed = $d.Content
$rngRed.Find.Execute("3D print must be RED, since the majority of the car needs to be red (it is in the rules)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$redStart = $rngRed.Start
$rngRed.Text = ""

$redInsert = $d.Range($redStart, $redStart)
$xmlRed = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>3D print must be RED, since the majority of the car needs to be red (it is in the rules)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$redInsert.InsertXML($xmlRed)

# ---------------------------------------------------------------------------
# 3) Remove the (now stale) <w:lastRenderedPageBreak/> that used to precede
#    "Preferable (" — the page now breaks earlier (at the paragraph above).
# ---------------------------------------------------------------------------
$rngPref = $d.Content
$rngPref.Find.Execute("Preferable (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prefStart = $rngPref.Start
$rngPref.Text = ""

$prefInsert = $d.Range($prefStart, $prefStart)
$xmlPref = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:bCs/></w:rPr><w:t>Preferable (</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$prefInsert.InsertXML($xmlPref)

Write-Output "edit applied"
